$d = $word.ActiveDocument

# The document contains 5 occurrences of a split "<id>...</id>" marker
# where the numeric id portion ("p093r_aN") was stored in its own run
# between the "<id>" and "</id>" runs. The edit collapses each of the
# three runs into a single run containing "<id>p093r_N</id>" (note the
# dropped "a" in the id), using the formatting of the first ("<id>") run.

for ($i = 1; $i -le 5; $i++) {
    $old = "<id>p093r_a$i</id>"
    $new = "<id>p093r_$i</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
}
